# Update the "date" column (F) values: shift each date forward by 2 days.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)  # Column F
    $cell.Value2 = $cell.Value2 + 2
}
